$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update B6, B7, B9 from name ("최준아") to "O" ---
$ws.Range("B6").Value = "O"
$ws.Range("B7").Value = "O"
$ws.Range("B9").Value = "O"

# --- Add "O" markers in column B for rows 23, 24, 26, 27 ---
$ws.Range("B23").Value = "O"
$ws.Range("B24").Value = "O"
$ws.Range("B26").Value = "O"
$ws.Range("B27").Value = "O"

# --- Add new rows 28-31 (stage-clear / ball & rabbit position tasks) ---
$ws.Range("A28").Value = "scrollview 위치 토끼위치로 변경"
$ws.Range("B28").Value = "O"

$ws.Range("A29").Value = "토끼 뛸 때 위치 anchoredposition `nset"
$ws.Range("A29").WrapText = $true
$ws.Rows.Item(29).RowHeight = 33
$ws.Range("B29").Value = "O"

$ws.Range("A30").Value = "구슬 위치 조정(max_y:4, min_y:0.85)`n구슬 시작시 위에서 아래로"
$ws.Range("A30").WrapText = $true
$ws.Rows.Item(30).RowHeight = 49.5
$ws.Range("B30").Value = "O"

$ws.Range("A31").Value = "위치 조정 시 shooter.possible 불가능하게.`n애니메이션작동할때나 떨어질때 구슬과 부딪히는 현상 없게"
$ws.Range("A31").WrapText = $true
$ws.Rows.Item(31).RowHeight = 66
$ws.Range("B31").Value = "O"

# --- Update sheet selection and scroll position ---
$ws.Range("B14").Select()

# --- Best-effort: restore window position (cosmetic, may be a no-op) ---
$win = $wb.Windows.Item(1)
$win.Left = 465
$win.Top = 165
